# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
# - Adds a new "Player Info" sheet (first tab) with player bio data.
# - On "ODI Batting": renames MATCH_CARD_LINK -> MATCH_CODE and rewrites
#   the column from a full scorecard URL to the bare numeric match code;
#   also clears the handful of stray empty INNING_NUMBER cells.
# - On "ODI Bowling": same MATCH_CARD_LINK -> MATCH_CODE rewrite.
# - Adds a new "ODI Batting Extra" sheet (last tab) with extra per-match
#   batting detail.

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $addr, $val) {
    # Force text storage so numeric-looking strings (match codes, "0",
    # "3.1", "9.39%", ...) don't silently get reinterpreted as numbers.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

function Strip-MatchCardLink($ws, $range) {
    $ws.Range($range).NumberFormat = "@"
    $ws.Range($range).Replace("http://www.howstat.com/cricket/Statistics/Matches/MatchScorecard_ODI.asp?MatchCode=", "")
}

# ---------------------------------------------------------------------
# 1. "Player Info" sheet - brand new, inserted before "ODI Batting"
# ---------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"
$playerInfo.Range("A1:D1").Font.Bold = $true

Set-TextCell $playerInfo "A2" "3986"
$playerInfo.Range("B2").Value = "James Douglas Sheahan Neesham"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium"

# ---------------------------------------------------------------------
# 2. "ODI Batting" - MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("D1").Value = "MATCH_CODE"
Strip-MatchCardLink $odiBatting "D2:D72"

# A handful of rows never had an INNING_NUMBER recorded (the player did
# not bat) and their B cell was a stray empty string - clear it outright
# so the cell no longer exists, matching the other "did not bat" rows.
$emptyInningRows = @(10, 11, 15, 35, 46, 48, 49, 51, 53, 65)
foreach ($r in $emptyInningRows) {
    $odiBatting.Range("B$r").ClearContents()
}

# ---------------------------------------------------------------------
# 3. "ODI Bowling" - MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------
$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$odiBowling.Range("B1").Value = "MATCH_CODE"
Strip-MatchCardLink $odiBowling "B2:B68"

# ---------------------------------------------------------------------
# 4. "ODI Batting Extra" sheet - brand new, appended at the end
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"
$extra.Range("A1:F1").Font.Bold = $true

Set-TextCell $extra "A2" "4315"
$extra.Range("B2").Value = 6
Set-TextCell $extra "F2" "YES"

Set-TextCell $extra "A3" "4328"
$extra.Range("B3").Value = 6
Set-TextCell $extra "C3" "3"
Set-TextCell $extra "D3" "0"
Set-TextCell $extra "E3" "9.39%"
Set-TextCell $extra "F3" "NO"

Set-TextCell $extra "A4" "4333"
Set-TextCell $extra "F4" "NO"

Set-TextCell $extra "A5" "4337"
Set-TextCell $extra "F5" "NO"

Set-TextCell $extra "A6" "4341"
$extra.Range("B6").Value = 7
Set-TextCell $extra "C6" "0"
Set-TextCell $extra "D6" "0"
Set-TextCell $extra "E6" "5.73%"
Set-TextCell $extra "F6" "NO"

Set-TextCell $extra "A7" "4346"
Set-TextCell $extra "F7" "NO"

Set-TextCell $extra "A8" "4353"
Set-TextCell $extra "F8" "NO"

Set-TextCell $extra "A9" "4355"
Set-TextCell $extra "F9" "NO"

Set-TextCell $extra "A10" "4402"
Set-TextCell $extra "F10" "NO"

Set-TextCell $extra "A11" "4406"
Set-TextCell $extra "F11" "NO"

Set-TextCell $extra "A12" "4410"
$extra.Range("B12").Value = 6
Set-TextCell $extra "C12" "1"
Set-TextCell $extra "D12" "1"
Set-TextCell $extra "E12" "6.33%"
Set-TextCell $extra "F12" "NO"

Set-TextCell $extra "A13" "4423"
$extra.Range("B13").Value = 6
Set-TextCell $extra "C13" "0"
Set-TextCell $extra "D13" "1"
Set-TextCell $extra "E13" "4.28%"
Set-TextCell $extra "F13" "NO"

Set-TextCell $extra "A14" "4452"
$extra.Range("B14").Value = 6
Set-TextCell $extra "F14" "NO"

Set-TextCell $extra "A15" "4453"
$extra.Range("B15").Value = 6
Set-TextCell $extra "C15" "2"
Set-TextCell $extra "D15" "0"
Set-TextCell $extra "E15" "10.91%"
Set-TextCell $extra "F15" "NO"

Set-TextCell $extra "A16" "4455"
$extra.Range("B16").Value = 7
Set-TextCell $extra "C16" "0"
Set-TextCell $extra "D16" "0"
Set-TextCell $extra "E16" "1.26%"
Set-TextCell $extra "F16" "NO"

Set-TextCell $extra "A17" "4639"
$extra.Range("B17").Value = 8
Set-TextCell $extra "C17" "0"
Set-TextCell $extra "D17" "0"
Set-TextCell $extra "E17" "0.47%"
Set-TextCell $extra "F17" "NO"

Set-TextCell $extra "A18" "4642"
Set-TextCell $extra "F18" "NO"

Set-TextCell $extra "A19" "4647"
Set-TextCell $extra "F19" "NO"

Set-TextCell $extra "A20" "4648"
$extra.Range("B20").Value = 7
Set-TextCell $extra "C20" "0"
Set-TextCell $extra "D20" "0"
Set-TextCell $extra "E20" "2.44%"
Set-TextCell $extra "F20" "NO"

Set-TextCell $extra "A21" "4649"
$extra.Range("B21").Value = 7
Set-TextCell $extra "C21" "5"
Set-TextCell $extra "D21" "0"
Set-TextCell $extra "E21" "14.88%"
Set-TextCell $extra "F21" "NO"
